$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '21.767.61'
$ws.Range("E2").Value = '  -1.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.540.85'
$ws.Range("E3").Value = '  -1.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.21'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3891'
$ws.Range("E7").Value = '  +2.53%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3190'
$ws.Range("E8").Value = '  -3.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.30'
$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07198'
$ws.Range("E10").Value = '  -2.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.063'
$ws.Range("E11").Value = '  -7.12%  '

$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.634'
$ws.Range("E13").Value = '  -3.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.56'
$ws.Range("E14").Value = '  -7.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.612'
$ws.Range("E15").Value = '  -4.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.543.73'
$ws.Range("E16").Value = '  -1.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001103'
$ws.Range("E17").Value = '  +0.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06594'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.36'
$ws.Range("E19").Value = '  -2.89%  '

$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.148'
$ws.Range("E21").Value = '  -5.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.40'
$ws.Range("E22").Value = '  -4.80%  '

$ws.Range("E23").Value = '  -7.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.362'
$ws.Range("E24").Value = '  +4.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '21.773.52'
$ws.Range("E25").Value = '  -1.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.388'
$ws.Range("E26").Value = '  -6.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '145.36'

$ws.Range("E28").Value = '  -3.81%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.847'
$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.717.77'
$ws.Range("E30").Value = '  -1.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.91'
$ws.Range("E31").Value = '  -3.07%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9646'
$ws.Range("E32").Value = '  -14.32%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.892'
$ws.Range("E33").Value = '  -2.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08216'
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.964'
$ws.Range("E35").Value = '  -4.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06128'
$ws.Range("E36").Value = '  -1.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.126'
$ws.Range("E37").Value = '  -3.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02211'
$ws.Range("E38").Value = '  -4.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.461'
$ws.Range("E39").Value = '  -21.90%  '

$ws.Range("E40").Value = '  -4.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -4.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5756'
$ws.Range("E44").Value = '  -4.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.10'
$ws.Range("E45").Value = '  -4.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.745'
$ws.Range("E46").Value = '  -0.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5519'
$ws.Range("E47").Value = '  -4.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.30'
$ws.Range("E48").Value = '  -2.44%  '

$ws.Range("E49").Value = '  -5.75%  '

$ws.Range("E50").Value = '  -3.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06736'
$ws.Range("E51").Value = '  -3.69%  '
